$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 40-43 (H7..H10), matching the formatting of row 39 ---
$ws.Range("A39:E39").Copy() | Out-Null
$ws.Range("A40:E43").PasteSpecial(-4122) | Out-Null

$ws.Range("A40").Value2 = "H7"
$ws.Range("B40").Value2 = 200
$ws.Range("C40").Value2 = 200
$ws.Range("D40").Value2 = 200
$ws.Range("E40").Value2 = 200

$ws.Range("A41").Value2 = "H8"
$ws.Range("B41").Value2 = 200
$ws.Range("C41").Value2 = 200
$ws.Range("D41").Value2 = 200
$ws.Range("E41").Value2 = 200

$ws.Range("A42").Value2 = "H9"
$ws.Range("B42").Value2 = 200
$ws.Range("C42").Value2 = 200
$ws.Range("D42").Value2 = 200
$ws.Range("E42").Value2 = 200

$ws.Range("A43").Value2 = "H10"
$ws.Range("B43").Value2 = 200
$ws.Range("C43").Value2 = 200
$ws.Range("D43").Value2 = 200
$ws.Range("E43").Value2 = 200

# --- View state: scrolled down a bit, selection moved ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("N17").Select() | Out-Null
